# Insert a new weekly price record for Cilantro at "Vega Monumental Concepción"
# right before the existing row 178 (2021-06-16 reading), pushing that row and
# every row after it down by one. The sheet is sorted by date, and the new
# record dated 2022-11-11 belongs between the rows currently at 177 and 178.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 178..240 down to 179..241, leaving a blank row 178 (formats
# inherited from the row above, matching Excel's normal insert behaviour).
$ws.Rows.Item(178).EntireRow.Insert()

# Populate the newly inserted row 178 with the new record.
$ws.Cells.Item(178, 1).Value = 11
$ws.Cells.Item(178, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(178, 3).Value = "Bíobío"
$ws.Cells.Item(178, 4).Value = 44876
$ws.Cells.Item(178, 5).Value = 8
$ws.Cells.Item(178, 6).Value = 100112040
$ws.Cells.Item(178, 7).Value = "Cilantro"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 220
$ws.Cells.Item(178, 11).Value = 19000
$ws.Cells.Item(178, 12).Value = 20000
$ws.Cells.Item(178, 13).Value = 19545
$ws.Cells.Item(178, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(178, 15).Value = "Región Metropolitana"
$ws.Cells.Item(178, 16).Value = 543
$ws.Cells.Item(178, 17).Value = 36
$ws.Cells.Item(178, 18).Value = "Hortaliza"
